$d = $word.ActiveDocument

function FindPos($needle) {
    $full = $d.Content.Text
    return $full.IndexOf($needle)
}

function MakeRun($start, $end) {
    $r = $d.Range($start, $end)
    $r.Font.Bold = 1
    $r.Font.Bold = 0
}

# =========================================================
# PASS 1 -- apply all text mutations (content only). Any
# run-merging side effects get cleaned up in PASS 2 below.
# =========================================================

# ---- Change 1: "...but to get to the llama." -> "...llama or the player." ----
$c1_needle = "but to get to the llama"
$c1_p = FindPos $c1_needle
$c1_p = $c1_p + $c1_needle.Length
$c1_text = " or the player"
$c1_ins = $d.Range($c1_p, $c1_p)
$c1_ins.InsertAfter($c1_text)

# ---- Change 2: rewrite the three-levels / ammunition paragraph span ----
$c2_old = "There will be a maximum of three levels with different variations of the maze and different placements of chickens or the enemies in this game. Level 1 will be the easiest of the three levels and level 2 will increase the size of the maze with a harder maze. Then level 3 will also increase the size of the maze even more, making it require more steps to get to the goal state of the game. In addition to the chickens and chest, there will also be a set of ammunition"
$c2_new = "There will be a maximum of three levels with different variations of the maze and different placements of chickens or the enemies in this game. Level 1 will be the easiest of the three levels with the least number of enemies. Level 2 will increase the number of enemies and will the maze more difficult to maneuver through. Lastly, level 3 will increase the difficulty in maneuvering through the maze. In addition to the chickens and the chest, there will also be a set of ammunition"
$c2_start = FindPos $c2_old
$c2_oldEnd = $c2_start + $c2_old.Length
$c2_range = $d.Range($c2_start, $c2_oldEnd)
$c2_range.Text = $c2_new

# ---- Change 3: "How to Play" paragraph ----
# 3a. "in this case is " -> "which is "
$c3a_old = "in this case is "
$c3a_new = "which is "
$c3a_start = FindPos $c3a_old
$c3a_oldEnd = $c3a_start + $c3a_old.Length
$c3a_range = $d.Range($c3a_start, $c3a_oldEnd)
$c3a_range.Text = $c3a_new

# 3b. move the _GoBack bookmark from the end of the paragraph to just
#     after "Once the game has started there will "
$c3b_needle = "Once the game has started there will "
$c3b_pos = FindPos $c3b_needle
$c3b_pos = $c3b_pos + $c3b_needle.Length

$goback = $d.Bookmarks.Item("_GoBack")
$goback.Delete()
$goback_range = $d.Range($c3b_pos, $c3b_pos)
$d.Bookmarks.Add("_GoBack", $goback_range)

# =========================================================
# PASS 2 -- re-split every run boundary the diff expects.
# Must run after all text/bookmark edits above because
# Range.Text assignment (and bookmark insertion) can merge
# neighbouring runs back together.
# =========================================================

# Change 1 split: "...llama" | " or the player" | ". ..."
$c1_p2 = $c1_p + $c1_text.Length
MakeRun $c1_p $c1_p2

# Change 2 splits
$c2_segs = @(
  "There will be a maximum of three levels with different variations of the maze and different placements of chickens or the enemies in this game. Level 1 will be the easiest of the",
  " three levels with the least number of enemies. Level",
  " 2 will increase the ",
  "number of enemies and will the maze more difficult to maneuver through",
  ". ",
  "Lastly,",
  " level 3 will ",
  "increase the difficulty in maneuvering through the maze. ",
  "In addition to the chickens and ",
  "the ",
  "chest, there will also be a set of ammunition that the player must get to"
)
$pos = $c2_start
foreach ($seg in $c2_segs) {
    $segEnd = $pos + $seg.Length
    MakeRun $pos $segEnd
    $pos = $segEnd
}

# Change 3 splits
# re-resolve "which is" position fresh (text already edited above)
$c3_which_needle = "you must toggle the shoot button, "
$c3_p1 = FindPos $c3_which_needle
$c3_p1 = $c3_p1 + $c3_which_needle.Length
$c3_whichis = "which is"
$c3_p2 = $c3_p1 + $c3_whichis.Length
MakeRun $c3_p1 $c3_p2

# "be a menu screen " run boundary (starts right after the bookmark point)
$c3_menu = "be a menu screen "
$c3_p3 = $c3b_pos + $c3_menu.Length
MakeRun $c3b_pos $c3_p3

Write-Output "All edits applied."
